$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a bad date value that had been typo'd in column A (45437 -> 45436)
$ws.Range("A3").Value2 = 45436

# Differentiate the pg data types onto their own columns: duplicate the
# date (numeric) values from column A into column B, and move the boolean
# values that used to live in column B over into column C.
foreach ($r in 1..3) {
    $colA = $ws.Cells.Item($r, 1)
    $colB = $ws.Cells.Item($r, 2)
    $colC = $ws.Cells.Item($r, 3)

    $colC.NumberFormat = $colB.NumberFormat
    $colB.NumberFormat = $colA.NumberFormat

    $colC.Value2 = $colB.Value2
    $colB.Value2 = $colA.Value2
}

$ws.Rows.Item(1).RowHeight = 13.8
$ws.Columns.Item(1).UseStandardWidth = $true

$ws.Range("B1:B3").Select()
